{"js": "// Replace the 25 \"two-digit \u00f7 one-digit =\" problems in the practice table\n// with the new values from the commit. Pairs are listed in document order\n// (top-to-bottom, left-to-right within each row). One value (\"44\u00f72=\")\n// occurs twice with two different replacements, so we resolve matches by\n// running one search per distinct \"from\" text and consuming its results\n// in document order.\nconst pairs = [\n  [\"71\u00f73=\", \"11\u00f78=\"],\n  [\"48\u00f79=\", \"76\u00f78=\"],\n  [\"90\u00f72=\", \"51\u00f76=\"],\n  [\"87\u00f74=\", \"89\u00f79=\"],\n  [\"35\u00f79=\", \"17\u00f76=\"],\n  [\"37\u00f76=\", \"30\u00f78=\"],\n  [\"23\u00f72=\", \"39\u00f79=\"],\n  [\"26\u00f73=\", \"67\u00f78=\"],\n  [\"14\u00f78=\", \"60\u00f78=\"],\n  [\"17\u00f75=\", \"48\u00f73=\"],\n  [\"44\u00f72=\", \"78\u00f77=\"],\n  [\"20\u00f78=\", \"61\u00f76=\"],\n  [\"37\u00f79=\", \"55\u00f77=\"],\n  [\"71\u00f76=\", \"52\u00f72=\"],\n  [\"44\u00f72=\", \"81\u00f76=\"],\n  [\"92\u00f74=\", \"58\u00f75=\"],\n  [\"82\u00f73=\", \"19\u00f78=\"],\n  [\"58\u00f79=\", \"92\u00f79=\"],\n  [\"68\u00f78=\", \"98\u00f79=\"],\n  [\"24\u00f72=\", \"59\u00f79=\"],\n  [\"59\u00f78=\", \"75\u00f77=\"],\n  [\"41\u00f77=\", \"55\u00f75=\"],\n  [\"96\u00f77=\", \"51\u00f75=\"],\n  [\"95\u00f72=\", \"33\u00f72=\"],\n  [\"54\u00f75=\", \"67\u00f78=\"],\n];\n\n// Run a search for every distinct \"from\" text used above.\nconst distinctFrom = [...new Set(pairs.map((p) => p[0]))];\nconst searchResults = {};\nfor (const from of distinctFrom) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  searchResults[from] = results;\n}\nawait context.sync();\n\n// Consume the search hits in document order, one per pair, in the order\n// the pairs are declared (which mirrors the order the diff lists them).\nconst nextIndex = {};\nfor (const [from, to] of pairs) {\n  const idx = nextIndex[from] || 0;\n  nextIndex[from] = idx + 1;\n  const range = searchResults[from].items[idx];\n  range.insertText(to, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($doc, [string]$findText, [string]$replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # MatchCase=True, Wrap=wdFindStop(1), Replace=wdReplaceOne(1): only the\n    # single next (first remaining) occurrence is substituted.\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1) | Out-Null\n}\n\nReplace-FirstMatch $d \"71\u00f73=\" \"11\u00f78=\"\nReplace-FirstMatch $d \"48\u00f79=\" \"76\u00f78=\"\nReplace-FirstMatch $d \"90\u00f72=\" \"51\u00f76=\"\nReplace-FirstMatch $d \"87\u00f74=\" \"89\u00f79=\"\nReplace-FirstMatch $d \"35\u00f79=\" \"17\u00f76=\"\nReplace-FirstMatch $d \"37\u00f76=\" \"30\u00f78=\"\nReplace-FirstMatch $d \"23\u00f72=\" \"39\u00f79=\"\nReplace-FirstMatch $d \"26\u00f73=\" \"67\u00f78=\"\nReplace-FirstMatch $d \"14\u00f78=\" \"60\u00f78=\"\nReplace-FirstMatch $d \"17\u00f75=\" \"48\u00f73=\"\nReplace-FirstMatch $d \"44\u00f72=\" \"78\u00f77=\"\nReplace-FirstMatch $d \"20\u00f78=\" \"61\u00f76=\"\nReplace-FirstMatch $d \"37\u00f79=\" \"55\u00f77=\"\nReplace-FirstMatch $d \"71\u00f76=\" \"52\u00f72=\"\nReplace-FirstMatch $d \"44\u00f72=\" \"81\u00f76=\"\nReplace-FirstMatch $d \"92\u00f74=\" \"58\u00f75=\"\nReplace-FirstMatch $d \"82\u00f73=\" \"19\u00f78=\"\nReplace-FirstMatch $d \"58\u00f79=\" \"92\u00f79=\"\nReplace-FirstMatch $d \"68\u00f78=\" \"98\u00f79=\"\nReplace-FirstMatch $d \"24\u00f72=\" \"59\u00f79=\"\nReplace-FirstMatch $d \"59\u00f78=\" \"75\u00f77=\"\nReplace-FirstMatch $d \"41\u00f77=\" \"55\u00f75=\"\nReplace-FirstMatch $d \"96\u00f77=\" \"51\u00f75=\"\nReplace-FirstMatch $d \"95\u00f72=\" \"33\u00f72=\"\nReplace-FirstMatch $d \"54\u00f75=\" \"67\u00f78=\"\n"}
